# Add a new "Spain" worksheet with the Zettler Spain market test data.
# The new sheet is built from a copy of the existing "Italy" sheet (same
# template: labels, styles, merged cells, page setup), then the two
# market-specific cells are updated to the Spain values.

$wb = $excel.ActiveWorkbook

$italy = $wb.Worksheets.Item("Italy")

# Duplicate the Italy sheet and place the copy right after it; Excel names
# the copy "Italy (2)" and makes it the active sheet/tab.
$italy.Copy($null, $italy)
$spain = $wb.Worksheets.Item("Italy (2)")
$spain.Name = "Spain"

# Market-specific values for the new sheet.
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3103/T2064"
$spain.Rows("3:4").RowHeight = 28.8

# Restore the Italy sheet's selection back to the full used range (it is no
# longer the active tab).
$italy.Select() | Out-Null
$italy.Range("A1:D21").Select() | Out-Null

# Leave the new Spain sheet selected, with B4 (the field we just edited)
# as the active cell - matching where the author's cursor ended up.
$spain.Select() | Out-Null
$spain.Range("B4").Select() | Out-Null
